# river update May 2024
#
# The "Chlorophyll A (92nd Percentile)" row (row 7) is dropped; the MCI and
# QMCI rows that followed it shift up to take rows 7 and 8 respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Delete()
